$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.449.95'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '3.825.73'
$ws.Range("E3").Value = '  +2.82%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''424.33'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").Value = '''130.61'
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("D7").Value = '3.817.67'
$ws.Range("E7").Value = '  +2.88%  '
$ws.Range("D8").Value = '''0.612'
$ws.Range("E8").Value = '  -5.36%  '
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '''0.730'
$ws.Range("E10").Value = '  -6.00%  '
$ws.Range("D11").Value = '''0.168'
$ws.Range("E11").Value = '  -9.55%  '
$ws.Range("D12").Value = '''0.0000366'
$ws.Range("E12").Value = '  -11.97%  '
$ws.Range("D13").Value = '''40.91'
$ws.Range("E13").Value = '  -5.39%  '
$ws.Range("D14").Value = '4.429.19'
$ws.Range("E14").Value = '  +3.30%  '
$ws.Range("D15").Value = '''10.07'
$ws.Range("E15").Value = '  -5.54%  '
$ws.Range("D16").Value = '''15.50'
$ws.Range("E16").Value = '  +15.63%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '''0.138'
$ws.Range("E17").Value = '  -1.37%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.827.31'
$ws.Range("E18").Value = '  +3.15%  '
$ws.Range("D19").Value = '''19.61'
$ws.Range("E19").Value = '  -5.92%  '
$ws.Range("D20").Value = '66.871.98'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '''1.07'
$ws.Range("E21").Value = '  -6.48%  '
$ws.Range("D22").Value = '''411.99'
$ws.Range("E22").Value = '  -8.14%  '
$ws.Range("D23").Value = '''14.45'
$ws.Range("E23").Value = '  -11.56%  '
$ws.Range("D24").Value = '''85.53'
$ws.Range("E24").Value = '  -5.19%  '
$ws.Range("D25").Value = '''3.04'
$ws.Range("E25").Value = '  -3.98%  '
$ws.Range("D26").Value = '''36.96'
$ws.Range("E26").Value = '  -2.18%  '
$ws.Range("D27").Value = '''5.68'
$ws.Range("E27").Value = '  +12.15%  '
$ws.Range("D28").Value = '''3.23'
$ws.Range("E28").Value = '  -3.18%  '
$ws.Range("D29").Value = '''9.50'
$ws.Range("D30").Value = '''689.54'
$ws.Range("E30").Value = '  +5.91%  '
$ws.Range("D31").Value = '''12.47'
$ws.Range("E31").Value = '  -2.32%  '
$ws.Range("E32").Value = '  -2.48%  '
$ws.Range("D33").Value = '''2.74'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("D34").Value = '''7.16'
$ws.Range("E34").Value = '  -2.05%  '
$ws.Range("D35").Value = '''0.152'
$ws.Range("E35").Value = '  -7.83%  '
$ws.Range("D36").Value = '''38.61'
$ws.Range("E36").Value = '  -8.30%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '0.0₃0795'
$ws.Range("E38").Value = '  +4.33%  '
$ws.Range("D39").Value = '''54.88'
$ws.Range("E39").Value = '  -4.36%  '
$ws.Range("D40").Value = '''3.15'
$ws.Range("E40").Value = '  +2.03%  '
$ws.Range("D41").Value = '''0.0456'
$ws.Range("E41").Value = '  -8.18%  '
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").Value = '''0.137'
$ws.Range("E43").Value = '  -8.57%  '
$ws.Range("D44").Value = '''148.98'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  +2.72%  '
$ws.Range("E46").Value = '  -3.20%  '
$ws.Range("E47").Value = '  -4.68%  '
$ws.Range("D48").Value = '''2.08'
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("D49").Value = '''26.17'
$ws.Range("E49").Value = '  -12.31%  '
$ws.Range("D50").Value = '''2.78'
$ws.Range("E50").Value = '  -4.38%  '
$ws.Range("D51").Value = '''2.54'
$ws.Range("E51").Value = '  -4.85%  '
